$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "Abalation-CAM: Visual Explanations for Deep Convolutional Network Via Gradient-free Localization [XAI-12]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/154"

$ws.Range("D32").Value = "HTTP 기본개념 - 1. 웹 브라우저가 메시지를 만든다."
$ws.Range("E32").Value = "https://dodonam.tistory.com/324"

$ws.Range("D51").Value = "[centos7] ""~은(는) sudoers 설정 파일에 없습니다"" 출력될 때 대처법"
$ws.Range("E51").Value = "https://bskyvision.com/1211"
